$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 240, pushing the existing rows 240..266 down to 241..267.
$ws.Rows.Item(240).Insert()

# Populate the newly inserted row 240 with the new weekly price record.
$ws.Range("A240").Value = 5
$ws.Range("B240").Value = "Macroferia Regional de Talca"
$ws.Range("C240").Value = "Maule"
$ws.Range("D240").Value = 44918
$ws.Range("E240").Value = 7
$ws.Range("F240").Value = 100112024
$ws.Range("G240").Value = "Choclo"
$ws.Range("H240").Value = "Choclero"
$ws.Range("I240").Value = "Primera"
$ws.Range("J240").Value = 20000
$ws.Range("K240").Value = 230
$ws.Range("L240").Value = 250
$ws.Range("M240").Value = 240
$ws.Range("N240").Value = "$/unidad"
$ws.Range("O240").Value = "Región del Maule"
$ws.Range("P240").Value = 240
$ws.Range("Q240").Value = 1
$ws.Range("R240").Value = "Hortaliza"

# Match the date formatting used by the other rows in column D.
$ws.Range("D240").NumberFormat = $ws.Range("D241").NumberFormat
